# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp (A1)
# - Refresh COVID-19 case/death counters for several countries (data rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated-as-of timestamp
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 22:03"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1206680
$ws.Range("C4").Value = 18558
$ws.Range("E4").Value = 955473
$ws.Range("F4").Value = 16043
$ws.Range("G4").Value = 708
$ws.Range("H4").Value = 69306

# Brasil (row 12)
$ws.Range("B12").Value = 105222
$ws.Range("C12").Value = 4075
$ws.Range("E12").Value = 54943
$ws.Range("G12").Value = 263
$ws.Range("H12").Value = 7288

# Suecia (row 25)
$ws.Range("D25").Value = 4074
$ws.Range("E25").Value = 15878

# Botsuana (row 184)
$ws.Range("C184").Value = 7
$ws.Range("D184").Value = 4
$ws.Range("E184").Value = 16
$ws.Range("G184").Value = 2
$ws.Range("H184").Value = 3

# Granada (row 185)
$ws.Range("B185").Value = 23
$ws.Range("D185").Value = 8
$ws.Range("E185").Value = 14
$ws.Range("F185").Value = 0
$ws.Range("H185").Value = 1

# Laos (row 186)
$ws.Range("B186").Value = 21
$ws.Range("D186").Value = 13
$ws.Range("E186").Value = 8
$ws.Range("F186").Value = 4

# Fiyi (row 187)
$ws.Range("B187").Value = 19
$ws.Range("D187").Value = 9
$ws.Range("E187").Value = 10

# Santa Lucia (row 188)
$ws.Range("D188").Value = 14
$ws.Range("E188").Value = 4

# Nueva Caledonia (row 190)
$ws.Range("D190").Value = 15
$ws.Range("E190").Value = 3
$ws.Range("F190").Value = 0

# Islas Virgenes de los Estados Unidos (row 191)
$ws.Range("B191").Value = 18
$ws.Range("D191").Value = 17
$ws.Range("E191").Value = 1
$ws.Range("F191").Value = 1

# San Vicente y las Granadinas (row 192)
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 0
$ws.Range("E192").Value = 17

# Gambia (row 193)
$ws.Range("C193").Value = 1
$ws.Range("E193").Value = 8
$ws.Range("H193").Value = 0

# Santo Tome y Principe (row 194)
$ws.Range("B194").Value = 17
$ws.Range("D194").Value = 9
$ws.Range("E194").Value = 7

# Curazao (row 198)
$ws.Range("D198").Value = 7
$ws.Range("H198").Value = 1

# San Cristobal y Nieves (row 199)
$ws.Range("D199").Value = 8
$ws.Range("H199").Value = 0
